# Auto-generated edit script: update cryptos Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.486.23"
Set-TextValue $ws.Range("E2") "  +0.27%  "
Set-TextValue $ws.Range("D3") "2.493.45"
Set-TextValue $ws.Range("E3") "  -0.55%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "569.18"
Set-TextValue $ws.Range("E5") "  -0.35%  "
Set-TextValue $ws.Range("D6") "165.92"
Set-TextValue $ws.Range("E6") "  +0.46%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("E8") "  -0.12%  "
Set-TextValue $ws.Range("E9") "  +0.90%  "
Set-TextValue $ws.Range("E10") "  -0.57%  "
Set-TextValue $ws.Range("E11") "  -2.37%  "
Set-TextValue $ws.Range("E12") "  -0.16%  "
Set-TextValue $ws.Range("D13") "2.951.51"
Set-TextValue $ws.Range("E13") "  -0.54%  "
Set-TextValue $ws.Range("D14") "69.352.55"
Set-TextValue $ws.Range("E14") "  +0.28%  "
Set-TextValue $ws.Range("D15") "0.0000174"
Set-TextValue $ws.Range("E15") "  +0.10%  "
Set-TextValue $ws.Range("D16") "24.14"
Set-TextValue $ws.Range("E16") "  -2.44%  "
Set-TextValue $ws.Range("D17") "2.492.20"
Set-TextValue $ws.Range("E17") "  -0.91%  "
Set-TextValue $ws.Range("E18") "  -0.48%  "
Set-TextValue $ws.Range("E19") "  -1.05%  "
Set-TextValue $ws.Range("D20") "352.72"
Set-TextValue $ws.Range("E20") "  +1.29%  "
Set-TextValue $ws.Range("E21") "  +0.30%  "
Set-TextValue $ws.Range("E22") "  -3.77%  "
Set-TextValue $ws.Range("D23") "0.999"
Set-TextValue $ws.Range("E23") "  -0.08%  "
Set-TextValue $ws.Range("D24") "69.33"
Set-TextValue $ws.Range("E24") "  -1.21%  "
Set-TextValue $ws.Range("E25") "  -2.65%  "
Set-TextValue $ws.Range("D26") "2.623.56"
Set-TextValue $ws.Range("E26") "  -1.17%  "
Set-TextValue $ws.Range("E27") "  -2.21%  "
Set-TextValue $ws.Range("E28") "  +0.29%  "
Set-TextValue $ws.Range("D29") "0.0₃0870"
Set-TextValue $ws.Range("E29") "  -1.14%  "
Set-TextValue $ws.Range("D30") "7.53"
Set-TextValue $ws.Range("E30") "  -3.45%  "
Set-TextValue $ws.Range("D31") "3.57"
Set-TextValue $ws.Range("E31") "  +138.79%  "
Set-TextValue $ws.Range("E32") "  -2.91%  "
Set-TextValue $ws.Range("D33") "438.32"
Set-TextValue $ws.Range("E33") "  -4.63%  "
Set-TextValue $ws.Range("E34") "  +0.05%  "
Set-TextValue $ws.Range("E35") "  -0.67%  "
Set-TextValue $ws.Range("E36") "  -2.30%  "
Set-TextValue $ws.Range("E37") "  -2.68%  "
Set-TextValue $ws.Range("D38") "19.06"
Set-TextValue $ws.Range("E38") "  +0.00%  "
Set-TextValue $ws.Range("E39") "  -1.44%  "
Set-TextValue $ws.Range("E40") "  +0.04%  "
Set-TextValue $ws.Range("D41") "0.314"
Set-TextValue $ws.Range("E41") "  -0.91%  "
Set-TextValue $ws.Range("D42") "4.58"
Set-TextValue $ws.Range("E42") "  -2.03%  "
Set-TextValue $ws.Range("E43") "  -1.73%  "
Set-TextValue $ws.Range("E44") "  -2.09%  "
Set-TextValue $ws.Range("E45") "  -3.60%  "
Set-TextValue $ws.Range("D46") "139.14"
Set-TextValue $ws.Range("E46") "  -1.64%  "
Set-TextValue $ws.Range("E47") "  -0.77%  "
Set-TextValue $ws.Range("D48") "0.504"
Set-TextValue $ws.Range("E48") "  -2.40%  "
Set-TextValue $ws.Range("E49") "  -0.97%  "
Set-TextValue $ws.Range("E50") "  -0.64%  "
Set-TextValue $ws.Range("D51") "0.0925"
Set-TextValue $ws.Range("E51") "  -0.31%  "
